$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.788.33'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '2.174.05'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.95%  '

$ws.Range("E6").Value = '  -2.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.97'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.04%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.27'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.35%  '

$ws.Range("E14").Value = '  -3.19%  '

$ws.Range("D15").Value = '2.499.31'
$ws.Range("E15").Value = '  -1.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.23%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.139.79'
$ws.Range("E18").Value = '  -2.57%  '

$ws.Range("D19").Value = '41.573.79'
$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000101'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("E25").Value = '  -2.77%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.70%  '

$ws.Range("E29").Value = '  -3.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0774'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.68%  '

$ws.Range("E35").Value = '  -8.29%  '

$ws.Range("E36").Value = '  -2.78%  '

$ws.Range("E37").Value = '  -0.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0301'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.32%  '

$ws.Range("E41").Value = '  -2.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.27%  '

$ws.Range("E44").Value = '  -2.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0969'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.84%  '

$ws.Range("E50").Value = '  -6.42%  '

$ws.Range("E51").Value = '  -2.41%  '
